$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(465, 44539, 8, 62, 409.078912641858)
    ,@(466, 44540, 13, 62, 409.078912641858)
    ,@(467, 44541, 6, 66, 435.471100554236)
    ,@(468, 44542, 9, 53, 349.6964898390077)
    ,@(469, 44543, 14, 59, 389.2847717075746)
    ,@(470, 44544, 6, 57, 376.0886777513856)
    ,@(471, 44545, 0, 56, 369.4906307732911)
    ,@(472, 44546, 7, 55, 362.8925837951966)
    ,@(473, 44547, 2, 44, 290.3140670361573)
    ,@(474, 44548, 7, 45, 296.9121140142518)
    ,@(475, 44550, 9, 45, 296.9121140142518)
    ,@(476, 44551, 8, 39, 257.3238321456849)
    ,@(477, 44552, 4, 37, 244.1277381894959)
    ,@(478, 44553, 4, 41, 270.5199261018739)
    ,@(479, 44554, 3, 37, 244.1277381894959)
    ,@(480, 44555, 6, 41, 270.5199261018739)
    ,@(481, 44556, 9, 43, 283.7160200580628)
    ,@(482, 44557, 20, 54, 356.2945368171021)
    ,@(483, 44558, 28, 74, 488.2554763789918)
    ,@(484, 44559, 14, 84, 554.2359461599367)
    ,@(485, 44560, 10, 90, 593.8242280285035)
    ,@(486, 44561, 33, 120, 791.765637371338)
    ,@(487, 44562, 15, 129, 851.1480601741883)
    ,@(488, 44563, 20, 140, 923.7265769332278)
    ,@(489, 44564, 28, 148, 976.5109527579837)
    ,@(490, 44565, 13, 133, 877.5402480865663)
    ,@(491, 44566, 34, 153, 1009.501187648456)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Range("A464").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("A$r").Value = $item[1]
    $ws.Range("B$r").Value = $item[2]
    $ws.Range("C$r").Value = $item[3]
    $ws.Range("D$r").Value = $item[4]
}

$ws.Range("A1").Select() | Out-Null
